$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: dates (forced to text via quote-prefix, then style reset to avoid a stray number format)
$ws.Cells.Item(2, 1).Value = "'2024-05-06"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(3, 1).Value = "'2024-05-07"
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(4, 1).Value = "'2024-05-08"
$ws.Cells.Item(4, 1).Style = "Normal"
$ws.Cells.Item(5, 1).Value = "'2024-05-09"
$ws.Cells.Item(5, 1).Style = "Normal"
$ws.Cells.Item(6, 1).Value = "'2024-05-10"
$ws.Cells.Item(6, 1).Style = "Normal"
$ws.Cells.Item(7, 1).Value = "'2024-05-13"
$ws.Cells.Item(7, 1).Style = "Normal"
$ws.Cells.Item(8, 1).Value = "'2024-05-14"
$ws.Cells.Item(8, 1).Style = "Normal"
$ws.Cells.Item(9, 1).Value = "'2024-05-15"
$ws.Cells.Item(9, 1).Style = "Normal"
$ws.Cells.Item(10, 1).Value = "'2024-05-16"
$ws.Cells.Item(10, 1).Style = "Normal"
$ws.Cells.Item(11, 1).Value = "'2024-05-17"
$ws.Cells.Item(11, 1).Style = "Normal"
$ws.Cells.Item(12, 1).Value = "'2024-05-27"
$ws.Cells.Item(12, 1).Style = "Normal"
$ws.Cells.Item(13, 1).Value = "'2024-05-28"
$ws.Cells.Item(13, 1).Style = "Normal"
$ws.Cells.Item(14, 1).Value = "'2024-05-29"
$ws.Cells.Item(14, 1).Style = "Normal"
$ws.Cells.Item(15, 1).Value = "'2024-05-30"
$ws.Cells.Item(15, 1).Style = "Normal"
$ws.Cells.Item(16, 1).Value = "'2024-05-31"
$ws.Cells.Item(16, 1).Style = "Normal"
$ws.Cells.Item(17, 1).Value = "'2024-05-20"
$ws.Cells.Item(17, 1).Style = "Normal"
$ws.Cells.Item(18, 1).Value = "'2024-05-21"
$ws.Cells.Item(18, 1).Style = "Normal"
$ws.Cells.Item(19, 1).Value = "'2024-05-22"
$ws.Cells.Item(19, 1).Style = "Normal"
$ws.Cells.Item(20, 1).Value = "'2024-05-23"
$ws.Cells.Item(20, 1).Style = "Normal"
$ws.Cells.Item(21, 1).Value = "'2024-05-24"
$ws.Cells.Item(21, 1).Style = "Normal"

# Column B: Vollkost (meat) menu text
$ws.Cells.Item(2, 2).Value = "Hähnchenstreifen `"Zürcher Art`" g,p`nButternudeln a,g,p,a-1,2"
$ws.Cells.Item(3, 2).Value = "Kasslerbraten kalt 1,2 `nRemoulade c,g,k,p,1,4,12 `nRöstkartoffeln 2 `nSalatbeilage c,g,k,p"
$ws.Cells.Item(4, 2).Value = "Tortellini (Fleisch) a,c,g,p,a-1,1,2 `nGemüse-Sahnesoße g,p"
$ws.Cells.Item(5, 2).Value = "Hühnerfrikassee g,p`nReis"
$ws.Cells.Item(6, 2).Value = "Ofenkartoffel `nDill-Heringshappen-Ragout d,g,p,1 `nSalatbeilage c,g,k,p"
$ws.Cells.Item(7, 2).Value = "Nudel-Gemüseauflauf m. Fleisch a,g,p,a-1,2 `nTomaten-Kräuterpestosoße 2"
$ws.Cells.Item(8, 2).Value = "Asia-Geschnetzeltes vom Huhn i,m,2,6,13 `nReis"
$ws.Cells.Item(9, 2).Value = "Hackfleischmasse - Frikadelle/Hackbraten (Rind) a,c,k,a-1`nJoghurt-Minz-Dip g,p`nGemüse Couscous a,i,a-1"
$ws.Cells.Item(10, 2).Value = "Hähnchenbrust paniert a,a-1`nBratensoße `nKartoffelpüree g,m,p,2 `nSommergemüse"
$ws.Cells.Item(11, 2).Value = "Fischfilet Piccata a,c,d,g,p,a-1`nFarfalle (Pasta/Nudeln) a,a-1,2 `nTomatisierte Zucchini"
$ws.Cells.Item(12, 2).Value = "Gemüseeintopf m. Kartoffel i`nRauchfleisch (50g) 1,2 `nBaguette a,a-1"
$ws.Cells.Item(13, 2).Value = "Gnocchi-Gemüse-Pfanne c,2 `nKäsesoße g,p,12 `nHähnchenbruststreifen (80g)"
$ws.Cells.Item(14, 2).Value = "Kichererbsen-Gemüseragout i`nRindfleisch `nCouscous a,i,a-1"
$ws.Cells.Item(15, 2).Value = "Lasagne Bolognese a,c,g,i,p,a1,2 `nFruchtige Tomatensoße a,a-5`nBohnensalat"
$ws.Cells.Item(16, 2).Value = "Matjesfilettopf `"HausfrauenArt`" c,d,g,k,p,1,3,4,12 `nRöstkartoffeln 2 `nSalatbeilage c,g,k,p"
$ws.Cells.Item(17, 2).Value = "Penne a,a-1,2 `nger. Putenbrust 1,2,16 `nGemüse-Sahnesoße g,p"
$ws.Cells.Item(18, 2).Value = "Putengyros 2 `nTzatziki g,p`nReis `nTomaten-Gurkensalat"
$ws.Cells.Item(19, 2).Value = "Gemüse-Quiche mit Räuchelachs a,c,d,g,p,a-1`nSalatbeilage c,g,k,p"
$ws.Cells.Item(20, 2).Value = "Chicken Nuggets a,a-1,16 `nAnanassauce 8 `nGemüsereis"
$ws.Cells.Item(21, 2).Value = "Mini Ofenkartoffeln 2 `nRäucherlachs-MeerrettichQuark d,g,m,p,2,13 `nSalatbeilage c,g,k,p"

# Column C: Vegetarisch menu text
$ws.Cells.Item(2, 3).Value = "Tofugeschnetzeltes `"Zürcher Art`" f,g,p`nButternudeln a,g,p,a-1,2"
$ws.Cells.Item(3, 3).Value = "Makkaroni a,c,a-1,2 `nZucchini-Paprika in Pestosoße g,p,2"
$ws.Cells.Item(4, 3).Value = "Tortellini (Vegetarisch) a,c,g,p,a-1,2 `nGemüse-Sahnesoße g,p"
$ws.Cells.Item(5, 3).Value = "Vegetarische Maultaschen a,c,g,i,p,a-1`nGeschmorte Butterzwiebeln g,p`nRohkost / Gemüsesticks"
$ws.Cells.Item(6, 3).Value = "Ofenkartoffel `nRucola-Dörrtomatencreme g,m,p,2 `nSalatbeilage c,g,k,p"
$ws.Cells.Item(7, 3).Value = "Nudel-Gemüseauflauf a,g,p,a1,2 `nTomaten-Kräuterpestosoße 2"
$ws.Cells.Item(8, 3).Value = "Asia-Geschnetzeltes mit Tofu f,i,m,2,6,13 `nReis"
$ws.Cells.Item(9, 3).Value = "Falafelbällchen a,a-1`nJoghurt-Minz-Dip g,p`nGemüse Couscous a,i,a-1"
$ws.Cells.Item(10, 3).Value = "Milchreis g,p`nheiße Kirschen"
$ws.Cells.Item(11, 3).Value = "Zucchini Piccata a,c,g,p,a-1`nTomaten-Basilikumsauce g,p,2 `nFarfalle (Pasta/Nudeln) a,a-1,2"
$ws.Cells.Item(12, 3).Value = "Gemüseeintopf m. Kartoffel i`nBaguette a,a-1"
$ws.Cells.Item(13, 3).Value = "Gnocchi-Gemüse-Pfanne c,2 `nKäsesoße g,p,12"
$ws.Cells.Item(14, 3).Value = "Kichererbsen-Gemüseragout i`nCouscous a,i,a-1"
$ws.Cells.Item(15, 3).Value = "Gemüselasagne a,c,g,p,a-1,2 `nFruchtige Tomatensoße a,a-5`nBohnensalat"
$ws.Cells.Item(16, 3).Value = "Dampfnudel a,c,g,p,a-1`nVanillesoße g,p,12"
$ws.Cells.Item(17, 3).Value = "Penne a,a-1,2 `ngeriebener Käse g,p`nGemüse-Sahnesoße g,p"
$ws.Cells.Item(18, 3).Value = "gebackener Hirtenkäse g,p`nTzatziki g,p`nReis `nTomaten-Gurkensalat"
$ws.Cells.Item(19, 3).Value = "Gemüse-Quiche a,c,g,p,a-1`nKräuter-Dip g,p`nSalatbeilage c,g,k,p"
$ws.Cells.Item(20, 3).Value = "gebackene Frühlingsrolle a,c,f,a-1`nAnanassauce 8 `nGemüsereis"
$ws.Cells.Item(21, 3).Value = "Käsespätzle a,c,g,p,a-1`nGeschmolzene Zwiebeln 2 `nSalatbeilage c,g,k,p"
